# Add the new "Env" worksheet at the end of the workbook (sheetId 3)
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Env"

# Column sizing (approximate Mac-Excel re-measured widths)
$newSheet.Columns.Item(1).ColumnWidth = 17
$newSheet.Columns.Item(2).ColumnWidth = 25.3

# --- Apply cell styles first (this fixes the order new style records are
#     created in the style table, independent of cell values) ---

# Rows 6-8: col A centered-horizontal only; col B centered-horizontal + text format
$newSheet.Range("A6:A8").HorizontalAlignment = -4108
$newSheet.Range("B6:B8").HorizontalAlignment = -4108
$newSheet.Range("B6:B8").NumberFormat = "@"

# Row 5: col B centered horizontal only + text format (reuses style created above)
$newSheet.Range("B5").HorizontalAlignment = -4108
$newSheet.Range("B5").NumberFormat = "@"

# --- Now fill in values (order controls shared-string table order) ---

$newSheet.Range("A2").Value = "ZoneID"
$newSheet.Range("A1").Value = "WorldID"
$newSheet.Range("A3").Value = "PlatformType"

$newSheet.Range("B1").Value = 1
$newSheet.Range("B2").Value = 2
$newSheet.Range("B3").Value = 3

$newSheet.Range("A5").Value = "赠送的道具"
$newSheet.Range("B4").Value = "测试"
$newSheet.Range("B5").Value = "1001:2,1002:5,1003:10"

$newSheet.Range("A6").Value = "基本信息"
$newSheet.Range("B6").Value = "110,标题,内容,true"

$newSheet.Range("A7").Value = "比例"
$newSheet.Range("A8").Value = "系数"

$newSheet.Range("B7").Value = "1.35"
$newSheet.Range("B8").Value = "15.246879"

$newSheet.Range("A4").Value = "名字"

# --- Apply remaining styles for rows 1-4 (existing style indices 1 and 3, so
#     applying after the values are set does not affect the new-style order) ---
$newSheet.Range("A1:A5").HorizontalAlignment = -4108
$newSheet.Range("A1:A5").VerticalAlignment = -4108

$newSheet.Range("B1:B4").HorizontalAlignment = -4108
$newSheet.Range("B1:B4").VerticalAlignment = -4108
$newSheet.Range("B1:B4").NumberFormat = "@"

# View state + activate as last (selected) tab
$newSheet.Range("C9").Select()
$newSheet.Activate()
